# Refresh the crypto price / 1h-volume columns (D/E) with the newly scraped values,
# and the TRON <-> Wrapped liquid staked Ether 2.0 rank swap (rows 13/14),
# as produced by the "Updated cryptos list ... with GitHub Actions" run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.977.65'
$ws.Range('E2').Value = '  +4.29%  '
$ws.Range('D3').Value = '3.486.50'
$ws.Range('E3').Value = '  +4.13%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'409.62"
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D6').Value = "'132.69"
$ws.Range('E6').Value = '  +18.09%  '
$ws.Range('D7').Value = '3.479.19'
$ws.Range('E7').Value = '  +4.18%  '
$ws.Range('D8').Value = "'0.605"
$ws.Range('E8').Value = '  +3.04%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = "'0.698"
$ws.Range('E10').Value = '  +10.05%  '
$ws.Range('D11').Value = "'0.130"
$ws.Range('E11').Value = '  +31.21%  '
$ws.Range('D12').Value = "'43.23"
$ws.Range('E12').Value = '  +7.92%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = "'0.142"
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '4.025.93'
$ws.Range('E14').Value = '  +3.81%  '
$ws.Range('D15').Value = "'8.75"
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').Value = "'20.26"
$ws.Range('E16').Value = '  +4.62%  '
$ws.Range('D17').Value = '3.471.21'
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('D18').Value = '62.879.76'
$ws.Range('E18').Value = '  +4.45%  '
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').Value = "'10.85"
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = "'0.0000141"
$ws.Range('E21').Value = '  +27.03%  '
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').Value = "'83.18"
$ws.Range('E23').Value = '  +9.85%  '
$ws.Range('D24').Value = "'13.15"
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').Value = "'314.73"
$ws.Range('E25').Value = '  +4.02%  '
$ws.Range('D26').Value = "'3.19"
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  +6.30%  '
$ws.Range('D28').Value = "'8.30"
$ws.Range('E28').Value = '  +3.45%  '
$ws.Range('D29').Value = "'7.76"
$ws.Range('E29').Value = '  +3.65%  '
$ws.Range('E30').Value = '  -0.81%  '
$ws.Range('E31').Value = '  -2.24%  '
$ws.Range('E32').Value = '  +3.67%  '
$ws.Range('D33').Value = "'11.92"
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').Value = "'2.64"
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').Value = "'42.97"
$ws.Range('E35').Value = '  +7.79%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = "'0.0498"
$ws.Range('E37').Value = '  -2.43%  '
$ws.Range('D38').Value = "'52.44"
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').Value = "'3.62"
$ws.Range('E39').Value = '  +5.99%  '
$ws.Range('D40').Value = "'0.997"
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').Value = "'3.02"
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('D42').Value = "'0.127"
$ws.Range('E42').Value = '  +2.77%  '
$ws.Range('E43').Value = '  +3.83%  '
$ws.Range('D44').Value = "'138.38"
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').Value = "'17.58"
$ws.Range('E45').Value = '  +3.71%  '
$ws.Range('D46').Value = "'4.01"
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('D47').Value = "'0.288"
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('D48').Value = "'2.25"
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('D49').Value = "'22.47"
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').Value = '2.209.68'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('D51').Value = '3.826.33'
$ws.Range('E51').Value = '  +4.03%  '
